$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-18: update only the ticker cells whose text actually changes.
# Cells that are blank both before and after are left untouched so their
# existing (empty) cell stays as it was.

$ws.Cells.Item(2, 2).Value = "NSE:ADFFOODS"
$ws.Cells.Item(2, 3).Value = "NSE:ADANIPOWER"
$ws.Cells.Item(2, 4).Value = "NSE:ADANIENT"
$ws.Cells.Item(2, 5).Value = "NSE:HAL"
$ws.Cells.Item(2, 6).Value = ""
$ws.Cells.Item(3, 2).Value = "NSE:AGI"
$ws.Cells.Item(3, 3).Value = "NSE:ALMONDZ"
$ws.Cells.Item(3, 4).Value = "NSE:DELHIVERY"
$ws.Cells.Item(3, 6).Value = ""
$ws.Cells.Item(4, 2).Value = "NSE:ARVINDFASN"
$ws.Cells.Item(4, 3).Value = "NSE:ANMOL"
$ws.Cells.Item(4, 4).Value = "NSE:KALYANKJIL"
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(5, 2).Value = "NSE:BBTC"
$ws.Cells.Item(5, 3).Value = "NSE:ATUL"
$ws.Cells.Item(5, 6).Value = ""
$ws.Cells.Item(6, 2).Value = "NSE:BHAGYANGR"
$ws.Cells.Item(6, 3).Value = "NSE:AXISILVER"
$ws.Cells.Item(6, 6).Value = ""
$ws.Cells.Item(7, 2).Value = "NSE:BORORENEW"
$ws.Cells.Item(7, 3).Value = "NSE:CONTROLPR"
$ws.Cells.Item(7, 6).Value = ""
$ws.Cells.Item(8, 2).Value = "NSE:BPL"
$ws.Cells.Item(8, 3).Value = "NSE:EXPLEOSOL"
$ws.Cells.Item(8, 6).Value = ""
$ws.Cells.Item(9, 2).Value = "NSE:BUTTERFLY"
$ws.Cells.Item(9, 3).Value = "NSE:GREENPOWER"
$ws.Cells.Item(10, 2).Value = "NSE:CHEMCON"
$ws.Cells.Item(10, 3).Value = "NSE:HDFCGOLD"
$ws.Cells.Item(11, 2).Value = "NSE:CHOLAHLDNG"
$ws.Cells.Item(11, 3).Value = "NSE:HDFCSILVER"
$ws.Cells.Item(12, 2).Value = "NSE:COFFEEDAY"
$ws.Cells.Item(12, 3).Value = "NSE:ITI"
$ws.Cells.Item(13, 2).Value = "NSE:DVL"
$ws.Cells.Item(13, 3).Value = "NSE:JPPOWER"
$ws.Cells.Item(14, 2).Value = "NSE:GAEL"
$ws.Cells.Item(14, 3).Value = "NSE:NETWEB"
$ws.Cells.Item(15, 2).Value = "NSE:GALLANTT"
$ws.Cells.Item(15, 3).Value = "NSE:RADICO"
$ws.Cells.Item(16, 2).Value = "NSE:GARFIBRES"
$ws.Cells.Item(16, 3).Value = "NSE:RTNPOWER"
$ws.Cells.Item(17, 2).Value = "NSE:HARRMALAYA"
$ws.Cells.Item(18, 2).Value = "NSE:HEUBACHIND"

# --- New rows 19-40: set the index number (col A, copying the bold/border
# style already used by the existing index cells) and the ticker in col B.
# Columns C:F have no ticker for any of these new rows.
$ws.Range("A2").Copy($ws.Range("A19:A40")) | Out-Null

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "NSE:HIMATSEIDE"
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "NSE:HLVLTD"
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "NSE:INFOBEAN"
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "NSE:JASH"
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "NSE:JINDRILL"
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "NSE:KALAMANDIR"
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "NSE:LUMAXIND"
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "NSE:MAGADSUGAR"
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "NSE:MAHEPC"
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "NSE:MANAKSIA"
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "NSE:MHRIL"
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "NSE:MINDACORP"
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "NSE:OLECTRA"
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = "NSE:OSWALGREEN"
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = "NSE:PARAGMILK"
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = "NSE:PGHH"
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = "NSE:PNCINFRA"
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = "NSE:RAMCOIND"
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = "NSE:RANEHOLDIN"
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = "NSE:RIIL"
$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).Value = "NSE:ROHLTD"
$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).Value = "NSE:RRKABEL"

Write-Host "Data updated: A2:F40"
